$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("CommonElements")
$ws3 = $wb.Worksheets.Item("ContactDetails_Elements")

# New rows 41-46 on the CommonElements sheet.
$ws2.Range("A41").Value = "alert_SuccessMsg"
$ws2.Range("B41").Value = "document.querySelector('div[id=""toast-container""] div[role=""alert""]')"

$ws2.Range("A42").Value = "alert_closeBtn"
$ws2.Range("B42").Value = "document.querySelector('div[id=""toast-container""] button')"

$ws2.Range("B43").Value = "document.querySelector('ion-icon[aria-label=""information circle outline""]')"
$ws2.Range("A43").Value = "help_button"

$ws2.Range("B44").Value = "document.querySelector('form pre')"
$ws2.Range("A44").Value = "help_description"

$ws2.Range("A45").Value = "specialChar_error"
$ws2.Range("B45").Value = "document.querySelector('[msg*=""ALPHANUMERIC""]')"

$ws2.Range("A46").Value = "mandatoryFillToastMsg"
$ws2.Range("B46").Value = "document.querySelector('ion-toast[role=""status""]').shadowRoot.querySelector('div[part=""message""]')"

# Update the recorded selection on ContactDetails_Elements, then return focus
# to CommonElements so it stays the active sheet/tab.
[void]$ws3.Select()
[void]$ws3.Range("B21").Select()

[void]$ws2.Select()
[void]$ws2.Range("B49").Select()
